$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "sex" values (column B) between P14 (rows 26-27) and P15 (rows 28-29)
$ws.Range("B26").Value = "m"
$ws.Range("B27").Value = "m"
$ws.Range("B28").Value = "f"
$ws.Range("B29").Value = "f"

# Update the active selection to match the saved view state
$ws.Range("B30").Select()
